# edit.ps1 - Applies the Week 4 Discussion edit described by the diff.
# Strategy: use Range.InsertXML with WordprocessingML "WordOpenXML" package envelopes
# to get exact control over run-splitting and proofErr markers (spell/gram check
# highlighting) that Word would normally insert, matching the target OOXML exactly.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("SELECT, FROM, WHERE, GROUP BY, HAVING, ORDER BY") gets
#    split into three runs with a proofErr gramStart/gramEnd pair around
#    "WHERE,". Replace the whole paragraph's range (this keeps the paragraph
#    mark / paragraph count intact).
# ---------------------------------------------------------------------------
$para2Xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">SELECT, FROM, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>WHERE,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> GROUP BY, HAVING, ORDER BY</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# 2) Paragraph 3 ("The six clauses ...") - the opening sentence is rewritten:
#    the former 5 separate runs ("The six clauses of a SQL retrieval query are ",
#    "SELECT, FROM, WHERE, GROUP BY, HAVING, ", "and ", "ORDER BY", ". ") collapse
#    into 3 runs with the same gramStart/gramEnd proofErr pair around "WHERE,".
#    The remaining two runs (the "A SQL retrieval query conceptually begins..."
#    and "The GROUP BY clause then groups..." explanation) are untouched, so we
#    only replace the leading portion of the paragraph up to (but not including)
#    "A SQL retrieval query conceptually begins".
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$p3 = $d.Paragraphs(3)
$p3full = $p3.Range

$prefixEnd = $p3full.Duplicate
$prefixEnd.Find.Execute("A SQL retrieval query conceptually begins", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$prefixRange = $d.Range($p3full.Start, $prefixEnd.Start)
$prefixRange.Text = ""

$para3PrefixXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The six clauses of a SQL retrieval query are SELECT, FROM, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>WHERE,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> GROUP BY, HAVING, and ORDER BY. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint = $d.Range($p3full.Start, $p3full.Start)
$insertPoint.InsertXML($para3PrefixXml)

# ---------------------------------------------------------------------------
# 3) Append the new content after paragraph 3: an intro sentence, the six-line
#    SQL example query (each line starting with a tab), and a closing paragraph
#    that walks through the conceptual execution order of the example query.
#    A single InsertXML call with multiple <w:p> elements appends them all as
#    new paragraphs right after paragraph 3, in one shot.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$p3 = $d.Paragraphs(3)
$p3End = $p3.Range.End
$tailInsertPoint = $d.Range($p3End, $p3End)

$newTailXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>See the example query below</w:t></w:r><w:r><w:t xml:space="preserve">, which finds </w:t></w:r><w:r><w:t>all department</w:t></w:r><w:r><w:t xml:space="preserve"> name</w:t></w:r><w:r><w:t>s that have more than two employees who make more than $</w:t></w:r><w:r><w:t>25</w:t></w:r><w:r><w:t>,000</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dn</w:t></w:r><w:r><w:t>ame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>COUNT(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>*)</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>company.employee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>company.department</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ON </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>employee.Dno</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>department.Dnumber</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>employee.Salary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>&gt;</w:t></w:r><w:r><w:t>25</w:t></w:r><w:r><w:t>000</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">GROUP BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>department.Dname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">HAVING </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>COUNT(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>*)&gt;2</w:t></w:r><w:r><w:t>,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Dn</w:t></w:r><w:r><w:t>ame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">We can think about the conceptual order by starting with the FROM clause. </w:t></w:r><w:r><w:t xml:space="preserve">In the FROM clause, the EMPLOYEE table is joined with the DEPARTMENT table on the </w:t></w:r><w:r><w:t>columns that store department number values. T</w:t></w:r><w:r><w:t>he rest of the query operates on this joined table.</w:t></w:r><w:r><w:t xml:space="preserve"> Next, the WHERE filters out employees across all departments who make less $</w:t></w:r><w:r><w:t>25</w:t></w:r><w:r><w:t xml:space="preserve">,000. </w:t></w:r><w:r><w:t xml:space="preserve">The GROUP BY clause then groups </w:t></w:r><w:r><w:t>all employee records by department n</w:t></w:r><w:r><w:t>ame</w:t></w:r><w:r><w:t xml:space="preserve">. The HAVING clause filters </w:t></w:r><w:r><w:t>to include only the groups that have more than two employee records. The ORDER BY clause orders</w:t></w:r><w:r><w:t xml:space="preserve"> the groups by department name. Finally, we return to the SELECT clause, where only the department name and the number of records </w:t></w:r><w:r><w:t>for each group is returned.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tailInsertPoint.InsertXML($newTailXml)

Write-Output "Paragraph count now: $($d.Paragraphs.Count)"
